$d = $word.ActiveDocument

# Replace the lone "13,5" occurrence (in the "- tickets  13,5" line) with "15"
$d.Content.Find.Execute("13,5", $true, $false, $false, $false, $false,
                         $true, 1, $false, "15", 2)
